$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44418
$ws.Cells.Item(2, 10).Value = 16
$ws.Cells.Item(2, 11).Value = 25000
$ws.Cells.Item(2, 12).Value = 26000
$ws.Cells.Item(2, 13).Value = 25500
$ws.Cells.Item(2, 16).Value = 1700

$ws.Cells.Item(3, 4).Value = 44351
$ws.Cells.Item(3, 10).Value = 34
$ws.Cells.Item(3, 11).Value = 24000
$ws.Cells.Item(3, 12).Value = 25000
$ws.Cells.Item(3, 13).Value = 24500
$ws.Cells.Item(3, 16).Value = 1633

$ws.Cells.Item(4, 4).Value = 44432
$ws.Cells.Item(4, 10).Value = 34
$ws.Cells.Item(4, 11).Value = 24000
$ws.Cells.Item(4, 12).Value = 25000
$ws.Cells.Item(4, 13).Value = 24500
$ws.Cells.Item(4, 16).Value = 1633

$ws.Cells.Item(5, 4).Value = 44385
$ws.Cells.Item(5, 10).Value = 25
$ws.Cells.Item(5, 11).Value = 14000
$ws.Cells.Item(5, 12).Value = 15000
$ws.Cells.Item(5, 13).Value = 14480
$ws.Cells.Item(5, 16).Value = 965

$ws.Cells.Item(6, 4).Value = 44446
$ws.Cells.Item(6, 10).Value = 34
$ws.Cells.Item(6, 11).Value = 24000
$ws.Cells.Item(6, 12).Value = 25000
$ws.Cells.Item(6, 13).Value = 24500
$ws.Cells.Item(6, 16).Value = 1633

$ws.Cells.Item(7, 4).Value = 44383
$ws.Cells.Item(7, 10).Value = 25
$ws.Cells.Item(7, 11).Value = 13000
$ws.Cells.Item(7, 12).Value = 14000
$ws.Cells.Item(7, 13).Value = 13480
$ws.Cells.Item(7, 16).Value = 899

$ws.Cells.Item(8, 4).Value = 44406
$ws.Cells.Item(8, 10).Value = 25
$ws.Cells.Item(8, 11).Value = 24000
$ws.Cells.Item(8, 12).Value = 25000
$ws.Cells.Item(8, 13).Value = 24520
$ws.Cells.Item(8, 16).Value = 1635

$ws.Cells.Item(9, 4).Value = 44425
$ws.Cells.Item(9, 10).Value = 25
$ws.Cells.Item(9, 11).Value = 24000
$ws.Cells.Item(9, 12).Value = 25000
$ws.Cells.Item(9, 13).Value = 24520
$ws.Cells.Item(9, 16).Value = 1635

$ws.Cells.Item(10, 4).Value = 44428
$ws.Cells.Item(10, 10).Value = 16
$ws.Cells.Item(10, 11).Value = 25000
$ws.Cells.Item(10, 12).Value = 26000
$ws.Cells.Item(10, 13).Value = 25500
$ws.Cells.Item(10, 16).Value = 1700

$ws.Cells.Item(11, 4).Value = 44449
$ws.Cells.Item(11, 10).Value = 18
$ws.Cells.Item(11, 11).Value = 24000
$ws.Cells.Item(11, 12).Value = 25000
$ws.Cells.Item(11, 13).Value = 24500
$ws.Cells.Item(11, 16).Value = 1633

$ws.Cells.Item(12, 4).Value = 44455
$ws.Cells.Item(12, 10).Value = 18
$ws.Cells.Item(12, 11).Value = 24000
$ws.Cells.Item(12, 12).Value = 25000
$ws.Cells.Item(12, 13).Value = 24500
$ws.Cells.Item(12, 16).Value = 1633

$ws.Cells.Item(13, 4).Value = 44397
$ws.Cells.Item(13, 10).Value = 34
$ws.Cells.Item(13, 11).Value = 23000
$ws.Cells.Item(13, 12).Value = 24000
$ws.Cells.Item(13, 13).Value = 23500
$ws.Cells.Item(13, 16).Value = 1567

$ws.Cells.Item(14, 4).Value = 44336
$ws.Cells.Item(14, 10).Value = 34
$ws.Cells.Item(14, 11).Value = 24000
$ws.Cells.Item(14, 12).Value = 25000
$ws.Cells.Item(14, 13).Value = 24500
$ws.Cells.Item(14, 16).Value = 1633

$ws.Cells.Item(15, 4).Value = 44442
$ws.Cells.Item(15, 10).Value = 28
$ws.Cells.Item(15, 11).Value = 24000
$ws.Cells.Item(15, 12).Value = 25000
$ws.Cells.Item(15, 13).Value = 24500
$ws.Cells.Item(15, 16).Value = 1633

$ws.Cells.Item(16, 4).Value = 44435
$ws.Cells.Item(16, 10).Value = 34
$ws.Cells.Item(16, 11).Value = 24000
$ws.Cells.Item(16, 12).Value = 25000
$ws.Cells.Item(16, 13).Value = 24500
$ws.Cells.Item(16, 16).Value = 1633

$ws.Cells.Item(17, 4).Value = 44349
$ws.Cells.Item(17, 10).Value = 21
$ws.Cells.Item(17, 11).Value = 24000
$ws.Cells.Item(17, 12).Value = 25000
$ws.Cells.Item(17, 13).Value = 24524
$ws.Cells.Item(17, 16).Value = 1635

$ws.Cells.Item(18, 4).Value = 44421
$ws.Cells.Item(18, 10).Value = 18
$ws.Cells.Item(18, 11).Value = 24000
$ws.Cells.Item(18, 12).Value = 25000
$ws.Cells.Item(18, 13).Value = 24500
$ws.Cells.Item(18, 16).Value = 1633

$ws.Cells.Item(19, 4).Value = 44400
$ws.Cells.Item(19, 10).Value = 16
$ws.Cells.Item(19, 11).Value = 24000
$ws.Cells.Item(19, 12).Value = 25000
$ws.Cells.Item(19, 13).Value = 24500
$ws.Cells.Item(19, 16).Value = 1633

$ws.Cells.Item(20, 4).Value = 44453
$ws.Cells.Item(20, 10).Value = 25
$ws.Cells.Item(20, 11).Value = 25000
$ws.Cells.Item(20, 12).Value = 26000
$ws.Cells.Item(20, 13).Value = 25520
$ws.Cells.Item(20, 16).Value = 1701

$ws.Cells.Item(21, 4).Value = 44413
$ws.Cells.Item(21, 10).Value = 25
$ws.Cells.Item(21, 11).Value = 24000
$ws.Cells.Item(21, 12).Value = 25000
$ws.Cells.Item(21, 13).Value = 24480
$ws.Cells.Item(21, 16).Value = 1632

$ws.Cells.Item(22, 4).Value = 44343
$ws.Cells.Item(22, 10).Value = 26
$ws.Cells.Item(22, 11).Value = 23000
$ws.Cells.Item(22, 12).Value = 24000
$ws.Cells.Item(22, 13).Value = 23500
$ws.Cells.Item(22, 16).Value = 1567

$ws.Cells.Item(23, 4).Value = 44329
$ws.Cells.Item(23, 10).Value = 25
$ws.Cells.Item(23, 11).Value = 23000
$ws.Cells.Item(23, 12).Value = 23000
$ws.Cells.Item(23, 13).Value = 23000
$ws.Cells.Item(23, 16).Value = 1533

$ws.Cells.Item(24, 4).Value = 44460
$ws.Cells.Item(24, 10).Value = 25
$ws.Cells.Item(24, 11).Value = 24000
$ws.Cells.Item(24, 12).Value = 25000
$ws.Cells.Item(24, 13).Value = 24480
$ws.Cells.Item(24, 16).Value = 1632

$ws.Cells.Item(25, 4).Value = 44411
$ws.Cells.Item(25, 10).Value = 34
$ws.Cells.Item(25, 11).Value = 25000
$ws.Cells.Item(25, 12).Value = 26000
$ws.Cells.Item(25, 13).Value = 25500
$ws.Cells.Item(25, 16).Value = 1700

$ws.Cells.Item(26, 4).Value = 44463
$ws.Cells.Item(26, 10).Value = 25
$ws.Cells.Item(26, 11).Value = 24000
$ws.Cells.Item(26, 12).Value = 25000
$ws.Cells.Item(26, 13).Value = 24480
$ws.Cells.Item(26, 16).Value = 1632

$ws.Cells.Item(27, 4).Value = 44341
$ws.Cells.Item(27, 10).Value = 36
$ws.Cells.Item(27, 11).Value = 24000
$ws.Cells.Item(27, 12).Value = 25000
$ws.Cells.Item(27, 13).Value = 24500
$ws.Cells.Item(27, 16).Value = 1633

$ws.Cells.Item(28, 4).Value = 44390
$ws.Cells.Item(28, 10).Value = 34
$ws.Cells.Item(28, 11).Value = 24000
$ws.Cells.Item(28, 12).Value = 25000
$ws.Cells.Item(28, 13).Value = 24500
$ws.Cells.Item(28, 16).Value = 1633

